$p = $ppt.ActivePresentation

# Remove the last three slides (slide10, slide11, slide12), which are
# no longer part of the deck. Delete from the end so indices stay valid.
$p.Slides.Item(12).Delete()
$p.Slides.Item(11).Delete()
$p.Slides.Item(10).Delete()
